# Courses.xlsx update: add the new 2020-Sem1-CAES-Wvl modules and their
# scheduled lab allocations, keeping the course list sorted alphabetically.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (unchanged, re-asserted for safety).
$headers = @("course", "labNum", "sectionNum", "allocatedDay", "allocatedTimeslot", "venueCapacity", "sessionLength")
for ($c = 0; $c -lt $headers.Length; $c++) {
    $ws.Cells.Item(1, $c + 1).Value = $headers[$c]
}

# Full allocation table: course, labNum, sectionNum, allocatedDay,
# allocatedTimeslot, venueCapacity, sessionLength.
$rows = @(
    @("BIOL101", 1, 1, "Monday", 2, 200, 1),
    @("BIOL101", 1, 2, "Tuesday", 2, 200, 1),
    @("BIOL101", 1, 3, "Wednesday", 2, 200, 1),
    @("BIOL101", 1, 4, "Thursday", 2, 200, 1),
    @("BIOL103", 1, 1, "Friday", 2, 200, 1),
    @("BIOL195", 1, 1, "Friday", 1, 200, 1),
    @("BIOL196", 1, 1, "Wednesday", 1, 200, 1),
    @("CHEM110", 1, 1, "Monday", 2, 144, 1),
    @("CHEM110", 1, 2, "Tuesday", 2, 144, 1),
    @("CHEM110", 1, 3, "Wednesday", 1, 144, 1),
    @("CHEM110", 1, 4, "Wednesday", 2, 144, 1),
    @("CHEM110", 1, 5, "Thursday", 2, 144, 1),
    @("CHEM110", 1, 6, "Friday", 2, 144, 1),
    @("CHEM195", 1, 1, "Tuesday", 1, 192, 1),
    @("CHEM195", 1, 2, "Wednesday", 2, 48, 1),
    @("CHEM196", 1, 1, "Monday", 2, 48, 1),
    @("CHEM196", 1, 2, "Tuesday", 2, 48, 1),
    @("CHEM196", 1, 3, "Thursday", 2, 48, 1),
    @("CHEM196", 1, 4, "Friday", 2, 48, 1),
    @("COMP100", 1, 1, "Monday", 2, 250, 1),
    @("COMP100", 1, 2, "Thursday", 2, 250, 1),
    @("GEOG110", 1, 1, "Wednesday", 2, 260, 1),
    @("GEOL101", 1, 1, "Monday", 2, 40, 1),
    @("GEOL101", 1, 2, "Thursday", 2, 40, 1),
    @("GEOL101", 1, 3, "Friday", 2, 40, 1),
    @("MATH130", 1, 1, "Tuesday", 2, 154, 1),
    @("MATH130", 1, 2, "Wednesday", 2, 136, 1),
    @("MATH134", 1, 1, "Wednesday", 1, 750, 1),
    @("MATH134", 1, 2, "Wednesday", 2, 500, 1),
    @("MATH134", 1, 3, "Thursday", 2, 250, 1),
    @("MATH140", 1, 1, "Friday", 2, 165, 1),
    @("MATH150", 1, 1, "Tuesday", 2, 491, 1),
    @("MATH150", 1, 2, "Wednesday", 1, 350, 1),
    @("MATH150", 1, 3, "Friday", 2, 491, 1),
    @("MATH196", 1, 1, "Wednesday", 1, 106, 1),
    @("PHYS113", 1, 1, "Wednesday", 2, 320, 1),
    @("PHYS131", 1, 1, "Monday", 2, 270, 1),
    @("PHYS131", 1, 2, "Tuesday", 2, 320, 1),
    @("PHYS131", 1, 3, "Wednesday", 1, 320, 1),
    @("PHYS131", 1, 4, "Thursday", 2, 320, 1),
    @("PHYS195", 1, 1, "Thursday", 1, 320, 1),
    @("STAT130", 1, 1, "Monday", 2, 491, 1),
    @("STAT130", 1, 2, "Wednesday", 1, 544, 1)
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2
    $row = $rows[$i]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
}
